# C5-PowerPoint.pptx edit
# 1) Slide 6's table switches from the custom "Table_0" table style to the
#    built-in PowerPoint table style {7DD8BBFC-82B9-4AA8-AB93-C09E2080FE4B}.
# 2) The deck's theme (ppt/theme/theme1.xml, used by the slide master) is
#    re-coloured from the "Integral" palette to the default Office palette
#    (same palette that, before this edit, only lived in the notes-master
#    theme part). Font / effect schemes were already identical between the
#    two theme parts, so only the 12 theme colours need to change.

$p = $ppt.ActivePresentation

# --- 1. Table style -------------------------------------------------------
$s = $p.Slides.Item(6)
$shp = $s.Shapes.Item(2)
$tbl = $shp.Table
$tbl.ApplyStyle("{7DD8BBFC-82B9-4AA8-AB93-C09E2080FE4B}")

# --- 2. Theme colours -------------------------------------------------------
$tcs = $p.SlideMaster.Theme.ThemeColorScheme

$tcs.Item(1).RGB  = 0         # dk1      000000
$tcs.Item(2).RGB  = 16777215  # lt1      FFFFFF
$tcs.Item(3).RGB  = 6968388   # dk2      44546A
$tcs.Item(4).RGB  = 15132391  # lt2      E7E6E6
$tcs.Item(5).RGB  = 13998939  # accent1  5B9BD5
$tcs.Item(6).RGB  = 3243501   # accent2  ED7D31
$tcs.Item(7).RGB  = 10855845  # accent3  A5A5A5
$tcs.Item(8).RGB  = 49407     # accent4  FFC000
$tcs.Item(9).RGB  = 12874308  # accent5  4472C4
$tcs.Item(10).RGB = 4697456   # accent6  70AD47
$tcs.Item(11).RGB = 12673797  # hlink    0563C1
$tcs.Item(12).RGB = 7491477   # folHlink 954F72
